$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Practice block (rows 2-5): fill in the pair_kind column (J) with "generic" ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- Trial rows (6-13): fill in the kind column (C) with "generic" ---
$ws.Range("C6").Value  = "generic"
$ws.Range("C7").Value  = "generic"
$ws.Range("C8").Value  = "generic"
$ws.Range("C9").Value  = "generic"
$ws.Range("C10").Value = "generic"
$ws.Range("C11").Value = "generic"
$ws.Range("C12").Value = "generic"
$ws.Range("C13").Value = "generic"

# --- New "stim details" block appended at the bottom of the sheet ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$stimRows = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$r = 29
foreach ($row in $stimRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
